# ---------------------------------------------------------------------------
# Applies the "testing redis cache" edit to measurements.xlsx:
#  - adds a third ("500 KB - 511 KB") data/column-group to the tables in
#    rows 4-6 (Nesting/Top level records vs response size)
#  - adds a new "gemidd x sneller" (avg-x-faster) summary row with formulas
#  - restyles/recolors the "No Cache" / "Redis Cache" result rows
#  - fixes a couple of data values and a label typo
#  - drops the now-redundant "Apollo Server" row and stray "ms"/note cells
#  - repositions the two screenshots that sit to the right of the table
# ---------------------------------------------------------------------------

function RGBV($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Cell content clean-up
# ---------------------------------------------------------------------------
$ws.Range("C1").ClearContents()          # old "ms" label next to the title
$ws.Range("O3").ClearContents()          # note text will move to column X

# ---------------------------------------------------------------------------
# 2. Row 4 header labels (KB ranges) + new "500 KB - 511 KB" column group
# ---------------------------------------------------------------------------
$ws.Range("F4").Value = "100 KB - 115 KB"     # typo fix: was "100 KB - 15 KB"
$ws.Range("J4").Value = "500 KB - 511 KB"     # new column group header
$ws.Range("J4:L4").Merge()

# ---------------------------------------------------------------------------
# 3. Row 5 / 6 data edits + new J:L column values
# ---------------------------------------------------------------------------
$ws.Range("H5").Value = 30                    # was 3
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 185
$ws.Range("L5").Value = 140

$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 3

# ---------------------------------------------------------------------------
# 4. New row 8 - "gemidd x sneller" (avg-x-faster) captions
# ---------------------------------------------------------------------------
$ws.Range("E8").Value = "gemidd x sneller"
$ws.Range("I8").Value = "gemidd x sneller"
$ws.Range("M8").Value = "gemidd x sneller"

# ---------------------------------------------------------------------------
# 5. Row 9 ("No Cache") - new columns + one corrected value
# ---------------------------------------------------------------------------
$ws.Range("H9").Value = 650                   # was 6300
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 430
$ws.Range("L9").Value = 2500

# ---------------------------------------------------------------------------
# 6. Row 10 ("Redis Cache") - corrected values + new columns + avg formulas
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = 6                    # was 4
$ws.Range("C10").Value = 12                   # was 16
$ws.Range("D10").Value = 24                   # was 70
$ws.Range("F10").Value = 18                   # was 16
$ws.Range("G10").Value = 25                   # was 90
$ws.Range("H10").Value = 50
$ws.Range("J10").Value = 63
$ws.Range("K10").Value = 65
$ws.Range("L10").Value = 200

$ws.Range("E10").Formula = "=SUM(B9:D9)/SUM(B10:D10)"
$ws.Range("I10").Formula = "=SUM(F9:H9)/SUM(F10:H10)"
$ws.Range("M10").Formula = "=SUM(J9:L9)/SUM(J10:L10)"

# ---------------------------------------------------------------------------
# 7. Drop the old "Apollo Server" row entirely
# ---------------------------------------------------------------------------
$ws.Range("A11").ClearContents()

# ---------------------------------------------------------------------------
# 8. A4 becomes a rich-text note: "Response Size" + bold caveat
# ---------------------------------------------------------------------------
$responseNote = "Response Size (in playground, niet 100% represatief met echte database data die gestored wordt in de cache)"
$ws.Range("A4").Value = $responseNote
$plainLen = ("Response Size").Length
$boldStart = $plainLen + 1
$boldLen = $responseNote.Length - $plainLen
$ws.Range("A4").Characters($boldStart, $boldLen).Font.Bold = $true

# ---------------------------------------------------------------------------
# 9. Move the explanatory note far right (column X) out of the table area
# ---------------------------------------------------------------------------
$ws.Range("X3").Value = "Afhankelijk van het aantal records opgevraagd van de nested element, veranderen de tijden"

# ---------------------------------------------------------------------------
# 10. Fill colours for the "No Cache" (gold) and "Redis Cache" (green) rows
# ---------------------------------------------------------------------------
$ws.Range("A9:M9").Interior.Color = RGBV 0xFF 0xF2 0xCC
$ws.Range("A10:M10").Interior.Color = RGBV 0xE2 0xF0 0xD9

# ---------------------------------------------------------------------------
# 11. Number format for the three average-speedup formula cells
# ---------------------------------------------------------------------------
$ws.Range("E10").Style = "Percent"
$ws.Range("E10").NumberFormat = "0.0"
$ws.Range("I10").NumberFormat = "0.0"
$ws.Range("M10").NumberFormat = "0.0"
$ws.Range("I10").Interior.Color = RGBV 0xE2 0xF0 0xD9
$ws.Range("M10").Interior.Color = RGBV 0xE2 0xF0 0xD9

# ---------------------------------------------------------------------------
# 12. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.28515625
$ws.Columns.Item(2).ColumnWidth = 8.7109375
$ws.Columns.Item(3).ColumnWidth = 8.7109375
$ws.Columns.Item(4).ColumnWidth = 10.140625
$ws.Columns.Item(5).ColumnWidth = 15.85546875
$ws.Columns.Item(9).ColumnWidth = 15.85546875
$ws.Columns.Item(13).ColumnWidth = 15.85546875

# ---------------------------------------------------------------------------
# 13. Selection cursor, as it ended up in the saved file
# ---------------------------------------------------------------------------
$ws.Range("E12").Select()

# ---------------------------------------------------------------------------
# 14. Reposition the two screenshots that float to the right of the table
#     (author dragged them further right once the table grew new columns)
# ---------------------------------------------------------------------------
$shapes = @($ws.Shapes)
$pic1 = $shapes[0]
$pic1.Left = 15544801 / 12700.0
$pic1.Top = 809624 / 12700.0
$pic1.Width = 1362074 / 12700.0
$pic1.Height = 2425155 / 12700.0

$pic2 = $shapes[1]
$pic2.Left = 16954501 / 12700.0
$pic2.Top = 828676 / 12700.0
$pic2.Width = 2091056 / 12700.0
$pic2.Height = 704850 / 12700.0
